$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume columns), matching the
# latest scrape. Column D (Price) values are forced to Text so that
# purely-numeric-looking strings (e.g. "168.20") keep their original
# formatting instead of being auto-coerced to a Double by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.883.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.068.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.654"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.371"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0779"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.884"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.363.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.085.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.875.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +5.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0222"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0969"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.319.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.248.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
